$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 56468
$ws.Range("D2").Value = 115550069
$ws.Range("C3").Value = 136920
$ws.Range("D3").Value = 322625031
$ws.Range("C4").Value = 49712
$ws.Range("D4").Value = 144940080
$ws.Range("C5").Value = 15758
$ws.Range("D5").Value = 53559040
$ws.Range("C6").Value = 5822
$ws.Range("D6").Value = 26603907
$ws.Range("C7").Value = 1174
$ws.Range("D7").Value = 6853421
$ws.Range("C8").Value = 75
$ws.Range("D8").Value = 521481
$ws.Range("C12").Value = 58820
$ws.Range("D12").Value = 93759141
$ws.Range("C13").Value = 14327
$ws.Range("D13").Value = 28987514
$ws.Range("C14").Value = 38498
$ws.Range("D14").Value = 89053212
$ws.Range("C15").Value = 12786
$ws.Range("D15").Value = 35503248
$ws.Range("C16").Value = 3685
$ws.Range("D16").Value = 11356675
$ws.Range("C17").Value = 1208
$ws.Range("D17").Value = 5217914
$ws.Range("C20").Value = 14367
$ws.Range("D20").Value = 22294426
$ws.Range("C21").Value = 20108
$ws.Range("D21").Value = 42625507
$ws.Range("C22").Value = 47634
$ws.Range("D22").Value = 114820550
$ws.Range("C23").Value = 16503
$ws.Range("D23").Value = 47761886
$ws.Range("C24").Value = 4896
$ws.Range("D24").Value = 15898972
$ws.Range("C25").Value = 1585
$ws.Range("D25").Value = 6689556
$ws.Range("C26").Value = 258
$ws.Range("D26").Value = 1414002
$ws.Range("C28").Value = 16088
$ws.Range("D28").Value = 24838140
$ws.Range("C29").Value = 11384
$ws.Range("D29").Value = 23325697
$ws.Range("C30").Value = 32927
$ws.Range("D30").Value = 74584031
$ws.Range("C31").Value = 11887
$ws.Range("D31").Value = 32281313
$ws.Range("C32").Value = 3273
$ws.Range("D32").Value = 9868612
$ws.Range("C33").Value = 1027
$ws.Range("D33").Value = 4356711
$ws.Range("C34").Value = 208
$ws.Range("D34").Value = 973490
$ws.Range("C36").Value = 11641
$ws.Range("D36").Value = 18093411
$ws.Range("C37").Value = 5041
$ws.Range("D37").Value = 11164450
$ws.Range("C38").Value = 11822
$ws.Range("D38").Value = 27585562
$ws.Range("C39").Value = 4894
$ws.Range("D39").Value = 14012335
$ws.Range("C40").Value = 1362
$ws.Range("D40").Value = 4482514
$ws.Range("C41").Value = 438
$ws.Range("D41").Value = 2099682
$ws.Range("C42").Value = 57
$ws.Range("D42").Value = 370438
$ws.Range("C44").Value = 3583
$ws.Range("D44").Value = 5522446
$ws.Range("C45").Value = 25816
$ws.Range("D45").Value = 53520901
$ws.Range("C46").Value = 76565
$ws.Range("D46").Value = 181077152
$ws.Range("C47").Value = 29470
$ws.Range("D47").Value = 82510586
$ws.Range("C48").Value = 9600
$ws.Range("D48").Value = 29401208
$ws.Range("C49").Value = 3315
$ws.Range("D49").Value = 13404524
$ws.Range("C50").Value = 585
$ws.Range("D50").Value = 3300708
$ws.Range("C51").Value = 30
$ws.Range("D51").Value = 142961
$ws.Range("C53").Value = 26603
$ws.Range("D53").Value = 49047371
$ws.Range("C54").Value = 2711
$ws.Range("D54").Value = 4393992
$ws.Range("C55").Value = 9055
$ws.Range("D55").Value = 15074433
$ws.Range("C56").Value = 3034
$ws.Range("D56").Value = 5349977
$ws.Range("C57").Value = 997
$ws.Range("D57").Value = 1971844
$ws.Range("C58").Value = 302
$ws.Range("D58").Value = 670437
$ws.Range("C61").Value = 9290
$ws.Range("D61").Value = 13848039
$ws.Range("C62").Value = 1829
$ws.Range("D62").Value = 3996094
$ws.Range("C63").Value = 4327
$ws.Range("D63").Value = 9411104
$ws.Range("C64").Value = 1732
$ws.Range("D64").Value = 3907584
$ws.Range("C68").Value = 2835
$ws.Range("D68").Value = 5609862
$ws.Range("C69").Value = 22970
$ws.Range("D69").Value = 45569721
$ws.Range("C70").Value = 66659
$ws.Range("D70").Value = 152658579
$ws.Range("C71").Value = 24390
$ws.Range("D71").Value = 68108753
$ws.Range("C72").Value = 7649
$ws.Range("D72").Value = 23329447
$ws.Range("C73").Value = 2477
$ws.Range("D73").Value = 10038509
$ws.Range("C74").Value = 519
$ws.Range("D74").Value = 2941713
$ws.Range("C75").Value = 28
$ws.Range("D75").Value = 98942
$ws.Range("C78").Value = 21313
$ws.Range("D78").Value = 32786319
$ws.Range("C79").Value = 83856
$ws.Range("D79").Value = 172913212
$ws.Range("C80").Value = 227756
$ws.Range("D80").Value = 515878363
$ws.Range("C81").Value = 102933
$ws.Range("D81").Value = 290380864
$ws.Range("C82").Value = 37577
$ws.Range("D82").Value = 127265720
$ws.Range("C83").Value = 13895
$ws.Range("D83").Value = 63119321
$ws.Range("C84").Value = 2747
$ws.Range("D84").Value = 17719212
$ws.Range("C85").Value = 154
$ws.Range("D85").Value = 890467
$ws.Range("C86").Value = 34
$ws.Range("D86").Value = 152125
$ws.Range("C90").Value = 79554
$ws.Range("D90").Value = 126626230
$ws.Range("C91").Value = 5649
$ws.Range("D91").Value = 8820971
$ws.Range("C92").Value = 13593
$ws.Range("D92").Value = 21543978
$ws.Range("C93").Value = 4360
$ws.Range("D93").Value = 7106451
$ws.Range("C95").Value = 422
$ws.Range("D95").Value = 957812
$ws.Range("C96").Value = 56
$ws.Range("D96").Value = 183975
$ws.Range("C98").Value = 6376
$ws.Range("D98").Value = 8779773
$ws.Range("C99").Value = 2313
$ws.Range("D99").Value = 4292019
$ws.Range("C100").Value = 7349
$ws.Range("D100").Value = 14653292
$ws.Range("C103").Value = 341
$ws.Range("D103").Value = 1100276
$ws.Range("C104").Value = 71
$ws.Range("D104").Value = 334891
$ws.Range("C106").Value = 4936
$ws.Range("D106").Value = 7229919
$ws.Range("C113").Value = 16451
$ws.Range("D113").Value = 34650842
$ws.Range("C114").Value = 43293
$ws.Range("D114").Value = 101628508
$ws.Range("C115").Value = 15355
$ws.Range("D115").Value = 43064220
$ws.Range("C116").Value = 4752
$ws.Range("D116").Value = 15112513
$ws.Range("C117").Value = 1478
$ws.Range("D117").Value = 6282980
$ws.Range("C122").Value = 13562
$ws.Range("D122").Value = 20859594
$ws.Range("C123").Value = 44835
$ws.Range("D123").Value = 91128972
$ws.Range("C124").Value = 95961
$ws.Range("D124").Value = 216652549
$ws.Range("C125").Value = 32436
$ws.Range("D125").Value = 87488479
$ws.Range("C126").Value = 10187
$ws.Range("D126").Value = 31136195
$ws.Range("C127").Value = 3241
$ws.Range("D127").Value = 13374589
$ws.Range("C128").Value = 670
$ws.Range("D128").Value = 3643934
$ws.Range("C132").Value = 35157
$ws.Range("D132").Value = 53845436
$ws.Range("C133").Value = 53952
$ws.Range("D133").Value = 110919973
$ws.Range("C134").Value = 113022
$ws.Range("D134").Value = 253243508
$ws.Range("C135").Value = 36597
$ws.Range("D135").Value = 101053850
$ws.Range("C136").Value = 10841
$ws.Range("D136").Value = 33593582
$ws.Range("C137").Value = 3435
$ws.Range("D137").Value = 14230729
$ws.Range("C138").Value = 568
$ws.Range("D138").Value = 3148927
$ws.Range("C139").Value = 50
$ws.Range("D139").Value = 245914
$ws.Range("C142").Value = 43919
$ws.Range("D142").Value = 65938128
$ws.Range("C143").Value = 19704
$ws.Range("D143").Value = 40577398
$ws.Range("C144").Value = 47880
$ws.Range("D144").Value = 112931415
$ws.Range("C145").Value = 17993
$ws.Range("D145").Value = 50542562
$ws.Range("C146").Value = 5189
$ws.Range("D146").Value = 16096507
$ws.Range("C147").Value = 1553
$ws.Range("D147").Value = 6653390
$ws.Range("C148").Value = 348
$ws.Range("D148").Value = 2005932
$ws.Range("C152").Value = 14843
$ws.Range("D152").Value = 23071971
$ws.Range("C153").Value = 53489
$ws.Range("D153").Value = 110989877
$ws.Range("C154").Value = 124246
$ws.Range("D154").Value = 287554965
$ws.Range("C155").Value = 39555
$ws.Range("D155").Value = 113970213
$ws.Range("C156").Value = 11827
$ws.Range("D156").Value = 39977043
$ws.Range("C157").Value = 4245
$ws.Range("D157").Value = 19032399
$ws.Range("C158").Value = 863
$ws.Range("D158").Value = 5272801
$ws.Range("C160").Value = 41224
$ws.Range("D160").Value = 64139481
